$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "34.950.88"
$cell.ClearFormats()
$ws.Range("E2").Value = "  -0.42%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.844.49"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +1.97%  "
$ws.Range("E4").Value = "  +0.07%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "232.06"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("E7").Value = "  +0.08%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "39.86"
$cell.ClearFormats()
$ws.Range("E8").Value = "  -1.66%  "
$ws.Range("E9").Value = "  +2.03%  "
$ws.Range("E10").Value = "  +0.51%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0981"
$cell.ClearFormats()
$ws.Range("E11").Value = "  -1.81%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "2.110.66"
$cell.ClearFormats()
$ws.Range("E12").Value = "  +1.91%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "11.58"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +4.61%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "1.842.26"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +1.77%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.677"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +1.93%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "4.65"
$cell.ClearFormats()
$ws.Range("E16").Value = "  -0.19%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "34.932.34"
$cell.ClearFormats()
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("E19").Value = "  -0.31%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "240.10"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +1.06%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "12.21"
$cell.ClearFormats()
$ws.Range("E21").Value = "  +2.19%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.69"
$cell.ClearFormats()
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("E24").Value = "  +2.02%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "171.96"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +0.12%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "7.82"
$cell.ClearFormats()
$ws.Range("E26").Value = "  -0.76%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "17.50"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +0.11%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "0.124"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +2.49%  "
$ws.Range("E29").Value = "  -2.62%  "
$ws.Range("E30").Value = "  +0.16%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.0553"
$cell.ClearFormats()
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("E32").Value = "  -3.72%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.97"
$cell.ClearFormats()
$ws.Range("E33").Value = "  -1.45%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.93"
$cell.ClearFormats()
$ws.Range("E34").Value = "  +9.54%  "
$ws.Range("E35").Value = "  +7.88%  "
$ws.Range("E36").Value = "  +16.14%  "
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("E38").Value = "  +7.42%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "90.52"
$cell.ClearFormats()
$ws.Range("E39").Value = "  -2.22%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.348.96"
$cell.ClearFormats()
$ws.Range("E40").Value = "  +2.58%  "
$ws.Range("E41").Value = "  +0.41%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "14.83"
$cell.ClearFormats()
$ws.Range("E42").Value = "  +3.01%  "
$ws.Range("E43").Value = "  +1.63%  "
$ws.Range("E44").Value = "  -2.27%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("E46").Value = "  +2.24%  "
$ws.Range("E47").Value = "  -0.91%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.027.84"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +2.00%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "3.41"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +21.00%  "
$ws.Range("E50").Value = "  +0.14%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.0670"
$cell.ClearFormats()
$ws.Range("E51").Value = "  +0.01%  "
